$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.784.66"
$ws.Range("E2").Value = "  +2.59%  "

$ws.Range("D3").Value = "1.696.47"
$ws.Range("E3").Value = "  +3.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.87"
$ws.Range("E5").Value = "  +2.92%  "

$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.10"
$ws.Range("E8").Value = "  +4.18%  "

$ws.Range("E9").Value = "  +2.37%  "

$ws.Range("E10").Value = "  +2.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0903"
$ws.Range("E11").Value = "  -1.53%  "

$ws.Range("D12").Value = "1.936.49"
$ws.Range("E12").Value = "  +3.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.75"
$ws.Range("E13").Value = "  +10.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.620"
$ws.Range("E14").Value = "  +6.99%  "

$ws.Range("D15").Value = "1.694.64"
$ws.Range("E15").Value = "  +3.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.03"
$ws.Range("E16").Value = "  +2.89%  "

$ws.Range("D17").Value = "30.787.52"
$ws.Range("E17").Value = "  +2.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.45"
$ws.Range("E18").Value = "  +2.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.72"
$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("D20").Value = "0.0₃0721"
$ws.Range("E20").Value = "  +1.60%  "

$ws.Range("E21").Value = "  -0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.24"
$ws.Range("E22").Value = "  +5.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.31"
$ws.Range("E23").Value = "  +2.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  +2.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.68"
$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.98"
$ws.Range("E26").Value = "  +1.59%  "

$ws.Range("E27").Value = "  +0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.77"
$ws.Range("E28").Value = "  +1.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("E30").Value = "  +2.10%  "

$ws.Range("E31").Value = "  +1.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.50"
$ws.Range("E32").Value = "  +3.13%  "

$ws.Range("D33").Value = "1.521.78"
$ws.Range("E33").Value = "  +6.22%  "

$ws.Range("E34").Value = "  +3.11%  "

$ws.Range("E35").Value = "  +5.05%  "

$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("E37").Value = "  +4.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "79.97"
$ws.Range("E38").Value = "  +6.89%  "

$ws.Range("E39").Value = "  -4.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.587"
$ws.Range("E40").Value = "  +5.09%  "

$ws.Range("E41").Value = "  +1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.858"
$ws.Range("E42").Value = "  +2.33%  "

$ws.Range("E43").Value = "  +1.43%  "

$ws.Range("E44").Value = "  +0.74%  "

$ws.Range("E45").Value = "  -1.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.57"
$ws.Range("E47").Value = "  -4.67%  "

$ws.Range("D48").Value = "1.828.76"
$ws.Range("E48").Value = "  +2.86%  "

$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.99"
$ws.Range("E50").Value = "  +6.07%  "

$ws.Range("D51").Value = "0.0₆0113"
$ws.Range("E51").Value = "  +2.43%  "
